# Apply the "Painkillers" + "Visibility" intent rows to the intents-en sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("intents-en")
$ws.Activate()

# Shared phrase text reused from earlier rows (yes/no training phrase sets).
$yesPhrases = "Yes|Okay I will|Why not|Yes that's alright|Yes I do|Exactly|Of course|Yep that's ok|Okay|Ok|I have it"
$noPhrases  = "No way|No|Na|I can't|No I cannot|Don't|Nope|I disagree|Of course not|No thanks|Not right now|Nah|I don't|Don't have it"

# Row 75: pain_killers_yes
$ws.Cells.Item(75,1).Value = "pain_killers_yes"
$ws.Cells.Item(75,2).Value = "en"
$ws.Cells.Item(75,3).Value = 1
$ws.Cells.Item(75,4).Value = 0
$ws.Cells.Item(75,7).Value = "PAINKILLER"
$ws.Cells.Item(75,8).Value = $yesPhrases
$ws.Rows.Item(75).RowHeight = 45

# Row 76: pain_killers_no
$ws.Cells.Item(76,1).Value = "pain_killers_no"
$ws.Cells.Item(76,2).Value = "en"
$ws.Cells.Item(76,3).Value = 1
$ws.Cells.Item(76,4).Value = 0
$ws.Cells.Item(76,7).Value = "PAINKILLERS"
$ws.Cells.Item(76,8).Value = $noPhrases
$ws.Rows.Item(76).RowHeight = 45

# Row 77: fallback_pain_killers
$ws.Cells.Item(77,1).Value = "fallback_pain_killers"
$ws.Cells.Item(77,2).Value = "en"
$ws.Cells.Item(77,3).Value = 0
$ws.Cells.Item(77,4).Value = 1
$ws.Cells.Item(77,5).Value = "PAINKILLERS"
$ws.Cells.Item(77,9).Value = "Again does pain killers work?"

# Row 78: visibility_set
$ws.Cells.Item(78,1).Value = "visibility_set"
$ws.Cells.Item(78,2).Value = "en"
$ws.Cells.Item(78,3).Value = 1
$ws.Cells.Item(78,4).Value = 0
$ws.Cells.Item(78,7).Value = "VISIBILITY"
$ws.Cells.Item(78,8).Value = "{@visibility:cloudy}|It is {@visibility:clear}|It looks {@visibility:murky}|I see it as {@visibility:clear}|very {@visibility:clear}|about {@visibility:clear}|around {@visibility:clear}"
$ws.Cells.Item(78,10).Value = "{@visibility:`$visibility:1:0}"
$ws.Rows.Item(78).RowHeight = 78.75

# Row 79: fallback_visibility
$ws.Cells.Item(79,1).Value = "fallback_visibility"
$ws.Cells.Item(79,2).Value = "en"
$ws.Cells.Item(79,3).Value = 0
$ws.Cells.Item(79,4).Value = 1
$ws.Cells.Item(79,5).Value = "VISIBILITY"
$ws.Cells.Item(79,9).Value = "Again how does it look?"

# Move the view roughly to where the new rows are, matching the saved selection.
$ws.Range("H78").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 70
$win.ScrollColumn = 1
